$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.223446
$ws.Range("H2").Value = 87.670338
$ws.Range("I2").Value = 0.0169041244192178
$ws.Range("J2").Value = 0.0169041244192178
$ws.Range("M2").Value = 35.04689966666667
$ws.Range("N2").Value = 105.140699
$ws.Range("O2").Value = 0.3824629895491901
$ws.Range("P2").Value = 0.3824629895491901
$ws.Range("Q2").Value = 1024.191179876251
$ws.Range("R2").Value = 9217.720618886264
$ws.Range("S2").Value = 0.006465201961085507
$ws.Range("T2").Value = 0.006465201961085507
$ws.Range("G3").Value = 29.223446
$ws.Range("H3").Value = 87.670338
$ws.Range("I3").Value = 0.0169041244192178
$ws.Range("J3").Value = 0.0169041244192178
$ws.Range("O3").Value = 0.3264402385872224
$ws.Range("P3").Value = 0.3264402385872223
$ws.Range("Q3").Value = 874.168801304974
$ws.Range("R3").Value = 7867.519211744766
$ws.Range("S3").Value = 0.005518186408517551
$ws.Range("T3").Value = 0.00551818640851755
$ws.Range("G4").Value = 29.223446
$ws.Range("H4").Value = 87.670338
$ws.Range("I4").Value = 0.0169041244192178
$ws.Range("J4").Value = 0.0169041244192178
$ws.Range("M4").Value = 8.911727666666666
$ws.Range("N4").Value = 26.735183
$ws.Range("O4").Value = 0.09725271102035077
$ws.Range("P4").Value = 0.09725271102035075
$ws.Range("Q4").Value = 260.4313922335393
$ws.Range("R4").Value = 2343.882530101854
$ws.Range("S4").Value = 0.001643971927194244
$ws.Range("T4").Value = 0.001643971927194243
$ws.Range("G5").Value = 29.223446
$ws.Range("H5").Value = 87.670338
$ws.Range("I5").Value = 0.0169041244192178
$ws.Range("J5").Value = 0.0169041244192178
$ws.Range("M5").Value = 17.76285166666667
$ws.Range("N5").Value = 53.288555
$ws.Range("O5").Value = 0.1938440608432367
$ws.Range("P5").Value = 0.1938440608432367
$ws.Range("Q5").Value = 519.0917364868433
$ws.Range("R5").Value = 4671.82562838159
$ws.Range("S5").Value = 0.0032767641224205
$ws.Range("T5").Value = 0.003276764122420499
$ws.Range("I6").Value = 0.9471112884046843
$ws.Range("J6").Value = 0.9471112884046842
$ws.Range("M6").Value = 35.04689966666667
$ws.Range("N6").Value = 105.140699
$ws.Range("O6").Value = 0.3824629895491901
$ws.Range("P6").Value = 0.3824629895491901
$ws.Range("Q6").Value = 57383.80787368788
$ws.Range("R6").Value = 516454.270863191
$ws.Range("S6").Value = 0.3622350147990407
$ws.Range("T6").Value = 0.3622350147990406
$ws.Range("I7").Value = 0.9471112884046843
$ws.Range("J7").Value = 0.9471112884046842
$ws.Range("O7").Value = 0.3264402385872224
$ws.Range("P7").Value = 0.3264402385872223
$ws.Range("S7").Value = 0.3091752349554767
$ws.Range("T7").Value = 0.3091752349554766
$ws.Range("I8").Value = 0.9471112884046843
$ws.Range("J8").Value = 0.9471112884046842
$ws.Range("M8").Value = 8.911727666666666
$ws.Range("N8").Value = 26.735183
$ws.Range("O8").Value = 0.09725271102035077
$ws.Range("P8").Value = 0.09725271102035075
$ws.Range("Q8").Value = 14591.55797261616
$ws.Range("R8").Value = 131324.0217535455
$ws.Range("S8").Value = 0.09210914043533286
$ws.Range("T8").Value = 0.09210914043533283
$ws.Range("I9").Value = 0.9471112884046843
$ws.Range("J9").Value = 0.9471112884046842
$ws.Range("M9").Value = 17.76285166666667
$ws.Range("N9").Value = 53.288555
$ws.Range("O9").Value = 0.1938440608432367
$ws.Range("P9").Value = 0.1938440608432367
$ws.Range("Q9").Value = 29083.88693503407
$ws.Range("R9").Value = 261754.9824153066
$ws.Range("S9").Value = 0.183591898214834
$ws.Range("T9").Value = 0.1835918982148339
$ws.Range("G10").Value = 37.39212666666667
$ws.Range("H10").Value = 112.17638
$ws.Range("I10").Value = 0.02162924801792661
$ws.Range("J10").Value = 0.0216292480179266
$ws.Range("M10").Value = 35.04689966666667
$ws.Range("N10").Value = 105.140699
$ws.Range("O10").Value = 0.3824629895491901
$ws.Range("P10").Value = 0.3824629895491901
$ws.Range("Q10").Value = 1310.478111609958
$ws.Range("R10").Value = 11794.30300448962
$ws.Range("S10").Value = 0.008272386858637105
$ws.Range("T10").Value = 0.008272386858637102
$ws.Range("G11").Value = 37.39212666666667
$ws.Range("H11").Value = 112.17638
$ws.Range("I11").Value = 0.02162924801792661
$ws.Range("J11").Value = 0.0216292480179266
$ws.Range("O11").Value = 0.3264402385872224
$ws.Range("P11").Value = 0.3264402385872223
$ws.Range("Q11").Value = 1118.520743462073
$ws.Range("R11").Value = 10066.68669115866
$ws.Range("S11").Value = 0.007060656883434169
$ws.Range("T11").Value = 0.007060656883434166
$ws.Range("G12").Value = 37.39212666666667
$ws.Range("H12").Value = 112.17638
$ws.Range("I12").Value = 0.02162924801792661
$ws.Range("J12").Value = 0.0216292480179266
$ws.Range("M12").Value = 8.911727666666666
$ws.Range("N12").Value = 26.735183
$ws.Range("O12").Value = 0.09725271102035077
$ws.Range("P12").Value = 0.09725271102035075
$ws.Range("Q12").Value = 333.2284497308378
$ws.Range("R12").Value = 2999.05604757754
$ws.Range("S12").Value = 0.002103503007074911
$ws.Range("T12").Value = 0.002103503007074911
$ws.Range("G13").Value = 37.39212666666667
$ws.Range("H13").Value = 112.17638
$ws.Range("I13").Value = 0.02162924801792661
$ws.Range("J13").Value = 0.0216292480179266
$ws.Range("M13").Value = 17.76285166666667
$ws.Range("N13").Value = 53.288555
$ws.Range("O13").Value = 0.1938440608432367
$ws.Range("P13").Value = 0.1938440608432367
$ws.Range("Q13").Value = 664.1907994812111
$ws.Range("R13").Value = 5977.717195330901
$ws.Range("S13").Value = 0.004192701268780423
$ws.Range("T13").Value = 0.004192701268780422
$ws.Range("G14").Value = 24.817167
$ws.Range("H14").Value = 74.45150100000001
$ws.Range("I14").Value = 0.01435533915817136
$ws.Range("J14").Value = 0.01435533915817136
$ws.Range("M14").Value = 35.04689966666667
$ws.Range("N14").Value = 105.140699
$ws.Range("O14").Value = 0.3824629895491901
$ws.Range("P14").Value = 0.3824629895491901
$ws.Range("Q14").Value = 869.7647618599111
$ws.Range("R14").Value = 7827.882856739201
$ws.Range("S14").Value = 0.005490385930426773
$ws.Range("T14").Value = 0.005490385930426772
$ws.Range("G15").Value = 24.817167
$ws.Range("H15").Value = 74.45150100000001
$ws.Range("I15").Value = 0.01435533915817136
$ws.Range("J15").Value = 0.01435533915817136
$ws.Range("O15").Value = 0.3264402385872224
$ws.Range("P15").Value = 0.3264402385872223
$ws.Range("Q15").Value = 742.362592288923
$ws.Range("R15").Value = 6681.263330600307
$ws.Range("S15").Value = 0.004686160339793955
$ws.Range("T15").Value = 0.004686160339793954
$ws.Range("G16").Value = 24.817167
$ws.Range("H16").Value = 74.45150100000001
$ws.Range("I16").Value = 0.01435533915817136
$ws.Range("J16").Value = 0.01435533915817136
$ws.Range("M16").Value = 8.911727666666666
$ws.Range("N16").Value = 26.735183
$ws.Range("O16").Value = 0.09725271102035077
$ws.Range("P16").Value = 0.09725271102035075
$ws.Range("Q16").Value = 221.163833762187
$ws.Range("R16").Value = 1990.474503859683
$ws.Range("S16").Value = 0.001396095650748765
$ws.Range("T16").Value = 0.001396095650748765
$ws.Range("G17").Value = 24.817167
$ws.Range("H17").Value = 74.45150100000001
$ws.Range("I17").Value = 0.01435533915817136
$ws.Range("J17").Value = 0.01435533915817136
$ws.Range("M17").Value = 17.76285166666667
$ws.Range("N17").Value = 53.288555
$ws.Range("O17").Value = 0.1938440608432367
$ws.Range("P17").Value = 0.1938440608432367
$ws.Range("Q17").Value = 440.823656207895
$ws.Range("R17").Value = 3967.412905871056
$ws.Range("S17").Value = 0.002782697237201868
$ws.Range("T17").Value = 0.002782697237201868
